# Update countries & provincias Spain
# - Reorder "Belice" / "Nueva Caledonia" rows (192/193) so that Nueva
#   Caledonia appears first (row 192) and Belice second (row 193), each
#   keeping its own data.
# - Refresh the "datos actualizados" timestamp in A1.
# - Update the statistic columns (B..H) for the rows whose figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update "last updated" timestamp (row 1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 19:39"

# ---------------------------------------------------------------------
# 2) Swap rows 192 (Belice) and 193 (Nueva Caledonia) - Nueva Caledonia
#    now comes first, keeping each country's own statistics attached.
# ---------------------------------------------------------------------
$a192 = $ws.Range("A192").Value2
$b192 = $ws.Range("B192").Value2
$c192 = $ws.Range("C192").Value2
$d192 = $ws.Range("D192").Value2
$e192 = $ws.Range("E192").Value2
$f192 = $ws.Range("F192").Value2
$g192 = $ws.Range("G192").Value2
$h192 = $ws.Range("H192").Value2

$a193 = $ws.Range("A193").Value2
$b193 = $ws.Range("B193").Value2
$c193 = $ws.Range("C193").Value2
$d193 = $ws.Range("D193").Value2
$e193 = $ws.Range("E193").Value2
$f193 = $ws.Range("F193").Value2
$g193 = $ws.Range("G193").Value2
$h193 = $ws.Range("H193").Value2

$ws.Range("A192").Value = $a193
$ws.Range("B192").Value = $b193
$ws.Range("C192").Value = $c193
$ws.Range("D192").Value = $d193
$ws.Range("E192").Value = $e193
$ws.Range("F192").Value = $f193
$ws.Range("G192").Value = $g193
$ws.Range("H192").Value = $h193

$ws.Range("A193").Value = $a192
$ws.Range("B193").Value = $b192
$ws.Range("C193").Value = $c192
$ws.Range("D193").Value = $d192
$ws.Range("E193").Value = $e192
$ws.Range("F193").Value = $f192
$ws.Range("G193").Value = $g192
$ws.Range("H193").Value = $h192

# ---------------------------------------------------------------------
# 3) Update statistic figures for the countries that changed
# ---------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1307430
$ws.Range("C4").Value = 14807
$ws.Range("D4").Value = 220554
$ws.Range("E4").Value = 1009025
$ws.Range("F4").Value = 16798
$ws.Range("G4").Value = 923
$ws.Range("H4").Value = 77851

# Row 10 - Alemania
$ws.Range("B10").Value = 169901
$ws.Range("C10").Value = 471
$ws.Range("E10").Value = 20797
$ws.Range("F10").Value = 1712
$ws.Range("G10").Value = 12
$ws.Range("H10").Value = 7404

# Row 16 - India
$ws.Range("B16").Value = 59642
$ws.Range("C16").Value = 3291
$ws.Range("D16").Value = 17883
$ws.Range("E16").Value = 39774
$ws.Range("G16").Value = 96
$ws.Range("H16").Value = 1985

# Row 28 - Irlanda
$ws.Range("B28").Value = 22541
$ws.Range("C28").Value = 156
$ws.Range("E28").Value = 4002
$ws.Range("G28").Value = 26
$ws.Range("H28").Value = 1429

# Row 33 - Israel
$ws.Range("B33").Value = 16436
$ws.Range("C33").Value = 55
$ws.Range("D33").Value = 11229
$ws.Range("E33").Value = 4962
$ws.Range("F33").Value = 78

# Row 57 - Argentina
$ws.Range("E57").Value = 3485
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = 285

# Row 59 - Kazajistan
$ws.Range("B59").Value = 4834
$ws.Range("C59").Value = 256
$ws.Range("D59").Value = 1631
$ws.Range("E59").Value = 3172

# Row 71 - Irak
$ws.Range("B71").Value = 2603
$ws.Range("C71").Value = 60
$ws.Range("D71").Value = 1661
$ws.Range("E71").Value = 838
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 104

# Row 72 - Uzbekistan
$ws.Range("B72").Value = 2325
$ws.Range("C72").Value = 27
$ws.Range("E72").Value = 540

# Row 80 - Islandia
$ws.Range("D80").Value = 1765
$ws.Range("E80").Value = 26

# Row 132 - Montenegro
$ws.Range("D132").Value = 267
$ws.Range("E132").Value = 49

# Row 148 - Suazilandia
$ws.Range("B148").Value = 159
$ws.Range("C148").Value = 6
$ws.Range("E148").Value = 145

# Row 182 - Zimbabue
$ws.Range("D182").Value = 9
$ws.Range("E182").Value = 21
